$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in grade values of 5 for row 29 (student #26), columns C-F
$ws.Range("C29").Value = 5
$ws.Range("D29").Value = 5
$ws.Range("E29").Value = 5
$ws.Range("F29").Value = 5

# Update the active selection to B4 (frozen pane bottomRight)
$ws.Range("B4").Select()
